{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change summary (per commit message / diff):\n//   1. Insert a comma after \"Mrs. Brown\" in \"Ich frag Mrs. Brown ob wir\n//      getrennt Nachsitzen k\u00f6nnen.\" -> \"... Mrs. Brown, ob wir ...\".\n//   2. Append a large new block of story text directly after\n//      \"... getrennt Nachsitzen k\u00f6nnen.\" at the end of the paragraph.\n\nconst body = context.document.body;\n\n// --- Step 1: \"Brown ob wi\" -> \"Brown, ob wi\" (insert comma) ---\nconst commaResults = body.search(\"Brown ob wi\", { matchCase: true });\ncommaResults.load(\"text\");\nawait context.sync();\n\nif (commaResults.items.length > 0) {\n  commaResults.items[0].insertText(\"Brown, ob wi\", \"Replace\");\n  await context.sync();\n}\n\n// --- Step 2: append the continuation of the story after \"k\u00f6nnen.\" ---\nconst endResults = body.search(\"getrennt Nachsitzen k\u00f6nnen.\", { matchCase: true });\nendResults.load(\"text\");\nawait context.sync();\n\nif (endResults.items.length > 0) {\n  const appendedText = \" Doch sie meint \u201eNein! Ihr seid in einer Klasse also m\u00fcsst ihr auch zusammenhalten. Deshalb werdet ihr jetzt nochmal extra l\u00e4nger nachsitzen. Keine Diskussion!\u201c Lisa schaut mich mit einem t\u00f6dlichen Blick an. Ich versuchte diesen so gut wie m\u00f6glich zu erwidern. Am Montag, auf dem Weg zu dem Geschichtsraum h\u00f6re ich Stimmen aus dem Klassenraum der 9c. Ich schleiche so leise wie m\u00f6glich n\u00e4her ran und h\u00f6re die Stimme meines Physiklehrers. Ich verstehe nicht alles, aber ich kann eindeutig h\u00f6ren das er \u00f6fter den Namen Charlie benutzt. Ich schaue auf meine Uhr und merke, dass ich schon 5 Minuten zu sp\u00e4t bin. So gerne ich auch weiter gelauscht h\u00e4tte entferne ich mich dem Klassenraum uns mache mich schleunigst auf den Weg zum Geschichtsraum 2. Dort angekommen sehe ich nur Mrs. Brown, Lisa konnte ich nirgends entdecken. Ich trete langsam ein und begr\u00fc\u00dfe Mrs. Brown mit einem normalen \u201eHallo\u201c das Mrs. Brown erwiderte. Nach einigen Minuten stille fragt mich Mrs. Brown, ob ich nicht zuf\u00e4llig Lisa gesehen h\u00e4tte. Das hatte ich nicht, aber ich wollte Lisa schlecht dastehen lassen und antworte mit einem \u201eJa, ich sah sie auf dem Weg hierher. Sie sagte, dass sie keine Lust auf Sie h\u00e4tte und lief davon.\u201c Mrs. Brown guckte \u00e4u\u00dferst angewidert aus dem Fenster, ich vermute das sie hofft Lisa irgendwo noch zu sehen. Der Rest das Nachsitzen bestand daraus einen Aufsatz \u00fcber das alte \u00c4gypten zu schreiben. Nach diesem langweiligen Nachmittag mach ich mich auf den weg nach Hause. Ich wohne in einer heruntergekommenen H\u00fctte. Meine Mutter ist vor einigen Tangen verstorben und mein Vater liegt mit einem Aneurysma auf dem Sofa. Wir haben nicht genug Geld, womit wir eine Operation kaufen k\u00f6nnten. Ich gehe in die K\u00fcche und mache eine Suppe f\u00fcr meinen Vater. Er bedankt sich mit einem nicken und f\u00e4ngt an die Suppe zu l\u00f6ffeln. Ich wollte gerade die Karotten zur\u00fcck in den K\u00fchlschrak bringen als meine Aufmerksamkeit von einem gro\u00dfen Zettel, der am K\u00fchlschrank h\u00e4ngt, angezogen wird. Ich betrachte den Zettel genauer und sehe das die Beerdigung meiner Mutter nicht stattfinden kann, da wir nicht gen\u00fcgen daf\u00fcr zahlen k\u00f6nnen. Niedergeschlagen und \u00e4u\u00dfert traurig lege ich die Karotten in den K\u00fchlschrank und packe meine Tasche. Ich verabschiede mich bei meinem Vater und mache mich auf den weg zur Arbeit. Vor dem Postamt sehe ich wie meine Kollegen mit meinem Chef \u00fcber irgendwas reden. Leise \u00f6ffne ich die T\u00fcr und sehe das auf meinem Schreibtisch ein gro\u00dfer Karton mit der Aufschrift \u201eF\u00fcr jeden Umzug der beste Transportkarton\u201c. Langsam drehe ich mich zu meine Chef um und bevor ich auch nur irgendetwas sagen beziehungsweise fragen kann nickt er und sagt: \u201eEs tut mir leid Mrs. Johnson aber leider m\u00fcssen wir sie hier feuern ich habe schon die Kartons und den Rest bereitgestellt.\u201c Ich kann es nicht fassen. Meine einzige M\u00f6glichkeit noch an Geld zu kommen ist einfach weg. Was soll ich den jetzt machen? Schlimmer kann der Tag nicht werden denk ich und gehe r\u00fcber zu meinem Schreibtisch, um zu packen. Gerade als ich den Weg nach Hause antrete sehe ich Lisa auf der anderen Seite der Stra\u00dfe. Auch wenn ich sie nicht mag, gehe ich r\u00fcber und frage sie, warum sie heute Nachmittag nicht beim Nachsitzen war. Lisa blickt zu Boden und sagt das sie geh\u00f6rt hat, wie ich Mrs. Brown gefragt habe, ob wir getrennt Nachsitzen k\u00f6nnen und sich darum dachte das es ihr egal sei das sie dann noch ein zweites Mal Nachsitzen m\u00fcsse aber sie wollte mir einen Gefallen tun.\";\n  endResults.items[0].insertText(appendedText, \"End\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document ($d below).\n#\n# Change summary (per commit message / diff):\n#   1. Insert a comma after \"Mrs. Brown\" in \"Ich frag Mrs. Brown ob wir\n#      getrennt Nachsitzen koennen.\" -> \"... Mrs. Brown, ob wir ...\".\n#   2. Append a large new block of story text directly after\n#      \"... getrennt Nachsitzen koennen.\" at the end of the paragraph.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: \"Brown ob wi\" -> \"Brown, ob wi\" (insert comma) ---\n$find1 = $d.Content.Find\n$find1.Text = \"Brown ob wi\"\n$find1.Replacement.Text = \"Brown, ob wi\"\n$find1.Execute([ref]\"Brown ob wi\", [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]\"Brown, ob wi\", [ref]2) | Out-Null\n\n# --- Step 2: append the continuation of the story after \"koennen.\" ---\n$endRange = $d.Content\n$find2 = $endRange.Find\n$find2.Text = \"getrennt Nachsitzen k\u00f6nnen.\"\n$found2 = $find2.Execute()\n\nif ($found2) {\n    $endRange.Collapse(0)\n    $appendedText = ' Doch sie meint \u201eNein! Ihr seid in einer Klasse also m\u00fcsst ihr auch zusammenhalten. Deshalb werdet ihr jetzt nochmal extra l\u00e4nger nachsitzen. Keine Diskussion!\u201c Lisa schaut mich mit einem t\u00f6dlichen Blick an. Ich versuchte diesen so gut wie m\u00f6glich zu erwidern. Am Montag, auf dem Weg zu dem Geschichtsraum h\u00f6re ich Stimmen aus dem Klassenraum der 9c. Ich schleiche so leise wie m\u00f6glich n\u00e4her ran und h\u00f6re die Stimme meines Physiklehrers. Ich verstehe nicht alles, aber ich kann eindeutig h\u00f6ren das er \u00f6fter den Namen Charlie benutzt. Ich schaue auf meine Uhr und merke, dass ich schon 5 Minuten zu sp\u00e4t bin. So gerne ich auch weiter gelauscht h\u00e4tte entferne ich mich dem Klassenraum uns mache mich schleunigst auf den Weg zum Geschichtsraum 2. Dort angekommen sehe ich nur Mrs. Brown, Lisa konnte ich nirgends entdecken. Ich trete langsam ein und begr\u00fc\u00dfe Mrs. Brown mit einem normalen \u201eHallo\u201c das Mrs. Brown erwiderte. Nach einigen Minuten stille fragt mich Mrs. Brown, ob ich nicht zuf\u00e4llig Lisa gesehen h\u00e4tte. Das hatte ich nicht, aber ich wollte Lisa schlecht dastehen lassen und antworte mit einem \u201eJa, ich sah sie auf dem Weg hierher. Sie sagte, dass sie keine Lust auf Sie h\u00e4tte und lief davon.\u201c Mrs. Brown guckte \u00e4u\u00dferst angewidert aus dem Fenster, ich vermute das sie hofft Lisa irgendwo noch zu sehen. Der Rest das Nachsitzen bestand daraus einen Aufsatz \u00fcber das alte \u00c4gypten zu schreiben. Nach diesem langweiligen Nachmittag mach ich mich auf den weg nach Hause. Ich wohne in einer heruntergekommenen H\u00fctte. Meine Mutter ist vor einigen Tangen verstorben und mein Vater liegt mit einem Aneurysma auf dem Sofa. Wir haben nicht genug Geld, womit wir eine Operation kaufen k\u00f6nnten. Ich gehe in die K\u00fcche und mache eine Suppe f\u00fcr meinen Vater. Er bedankt sich mit einem nicken und f\u00e4ngt an die Suppe zu l\u00f6ffeln. Ich wollte gerade die Karotten zur\u00fcck in den K\u00fchlschrak bringen als meine Aufmerksamkeit von einem gro\u00dfen Zettel, der am K\u00fchlschrank h\u00e4ngt, angezogen wird. Ich betrachte den Zettel genauer und sehe das die Beerdigung meiner Mutter nicht stattfinden kann, da wir nicht gen\u00fcgen daf\u00fcr zahlen k\u00f6nnen. Niedergeschlagen und \u00e4u\u00dfert traurig lege ich die Karotten in den K\u00fchlschrank und packe meine Tasche. Ich verabschiede mich bei meinem Vater und mache mich auf den weg zur Arbeit. Vor dem Postamt sehe ich wie meine Kollegen mit meinem Chef \u00fcber irgendwas reden. Leise \u00f6ffne ich die T\u00fcr und sehe das auf meinem Schreibtisch ein gro\u00dfer Karton mit der Aufschrift \u201eF\u00fcr jeden Umzug der beste Transportkarton\u201c. Langsam drehe ich mich zu meine Chef um und bevor ich auch nur irgendetwas sagen beziehungsweise fragen kann nickt er und sagt: \u201eEs tut mir leid Mrs. Johnson aber leider m\u00fcssen wir sie hier feuern ich habe schon die Kartons und den Rest bereitgestellt.\u201c Ich kann es nicht fassen. Meine einzige M\u00f6glichkeit noch an Geld zu kommen ist einfach weg. Was soll ich den jetzt machen? Schlimmer kann der Tag nicht werden denk ich und gehe r\u00fcber zu meinem Schreibtisch, um zu packen. Gerade als ich den Weg nach Hause antrete sehe ich Lisa auf der anderen Seite der Stra\u00dfe. Auch wenn ich sie nicht mag, gehe ich r\u00fcber und frage sie, warum sie heute Nachmittag nicht beim Nachsitzen war. Lisa blickt zu Boden und sagt das sie geh\u00f6rt hat, wie ich Mrs. Brown gefragt habe, ob wir getrennt Nachsitzen k\u00f6nnen und sich darum dachte das es ihr egal sei das sie dann noch ein zweites Mal Nachsitzen m\u00fcsse aber sie wollte mir einen Gefallen tun.'\n    $endRange.InsertAfter($appendedText)\n}\n\n$d.Saved = $false\n"}
